# Apply cell value updates from the crypto price refresh.
# Cells whose new text looks like a pure number need NumberFormat
# forced to Text ("@") before assignment so Excel keeps them as
# strings (matching the source data, which stores formatted price
# strings like "0.9994" or "1.000" as text, not numbers). The
# NumberFormat is reset back to General afterwards so styling is
# unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells receiving plain-number-looking strings
# so they stay text cells instead of being parsed as numbers.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Write the new values.
$ws.Range("D2").Value = "29.395.76"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.848.42"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "240.18"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("D6").Value = "0.6299"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.07634"
$ws.Range("E8").Value = "  +0.91%  "
$ws.Range("D9").Value = "0.2930"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("D10").Value = "24.55"
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("D11").Value = "0.07741"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.857.97"
$ws.Range("E12").Value = "  -6.40%  "
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "0.00001119"
$ws.Range("E13").Value = "  +12.37%  "
$ws.Range("D14").Value = "5.004"
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "83.77"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").Value = "2.112.93"
$ws.Range("E17").Value = "  -6.69%  "
$ws.Range("D18").Value = "6.177"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("D19").Value = "29.413.06"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "229.12"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "7.493"
$ws.Range("E23").Value = "  -0.74%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "157.36"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("D26").Value = "0.1397"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "8.348"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").Value = "1.468"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("D30").Value = "1.299"
$ws.Range("E30").Value = "  +3.55%  "
$ws.Range("E31").Value = "  -1.83%  "
$ws.Range("D32").Value = "4.117"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").Value = "4.032"
$ws.Range("E33").Value = "  +0.26%  "
$ws.Range("D34").Value = "1.851"
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").Value = "0.7112"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("D37").Value = "2.585"
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("D38").Value = "1.240.40"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.781"
$ws.Range("E39").Value = "  -0.81%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.01806"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").Value = "6.401"
$ws.Range("E41").Value = "  +5.05%  "
$ws.Range("D42").Value = "0.9039"
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "101.91"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "66.02"
$ws.Range("E45").Value = "  -0.24%  "
$ws.Range("D46").Value = "7.157"
$ws.Range("D47").Value = "0.00000000117"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "9.045"
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "1.686"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("D51").Value = "0.1120"
$ws.Range("E51").Value = "  -0.39%  "

# Restore General number format on the cells we forced to Text.
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D51").NumberFormat = "General"
